# Updated Bus 11 AM
# The row for "Bjørn" / "Lindi" (row 9) is removed from the attendee list.
# Deleting the entire row shifts all subsequent rows up by one and causes
# the two now-unused shared strings ("Bjørn", "Lindi") to be dropped from
# the shared string table on save.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the row for Bjørn Lindi (row 9), shifting rows 10-24 up to 9-23.
$ws.Rows.Item(9).Delete()

# Reproduce the post-delete selection state (whole row 9 selected).
$ws.Range("A9:XFD9").Select()
